$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3261.6042
$ws.Range("I137").Value = 3377
$ws.Range("J137").Value = 3007.7334
$ws.Range("K137").Value = 10131
$ws.Range("L137").Value = 9023.200199999999
$ws.Range("M137").Value = -7581
$ws.Range("N137").Value = -14123.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6284.82
$ws.Range("I32").Value = 5731.977
$ws.Range("J32").Value = 9984.615
$ws.Range("K32").Value = 5731.977
$ws.Range("L32").Value = 9984.615
$ws.Range("M32").Value = -5444.977
$ws.Range("N32").Value = -10558.615

$ws.Range("H45").Value = 1600.48
$ws.Range("I45").Value = 1250.5
$ws.Range("K45").Value = 1250.5
$ws.Range("M45").Value = -873.5

$ws.Range("H74").Value = 1641.279
$ws.Range("I74").Value = 1273.1082
$ws.Range("J74").Value = 3911.6667
$ws.Range("K74").Value = 1273.1082
$ws.Range("L74").Value = 3911.6667
$ws.Range("M74").Value = -399.1081999999999
$ws.Range("N74").Value = -5659.6667

$ws.Range("H77").Value = 1641.279
$ws.Range("I77").Value = 1273.1082
$ws.Range("J77").Value = 3911.6667
$ws.Range("K77").Value = 6365.540999999999
$ws.Range("L77").Value = 19558.3335
$ws.Range("M77").Value = -1997.540999999999
$ws.Range("N77").Value = -28294.3335

$ws.Range("H97").Value = 696.1539
$ws.Range("I97").Value = 654.1667
$ws.Range("K97").Value = 654.1667
$ws.Range("M97").Value = -158.1667

$ws.Range("H122").Value = 2338.4243
$ws.Range("I122").Value = 1564.2174
$ws.Range("K122").Value = 4692.6522
$ws.Range("M122").Value = -2242.6522

$ws.Range("H132").Value = 2061.2
$ws.Range("I132").Value = 1454.4117
$ws.Range("J132").Value = 5499.6665
$ws.Range("K132").Value = 4363.2351
$ws.Range("L132").Value = 16498.9995
$ws.Range("M132").Value = -1833.2351
$ws.Range("N132").Value = -21558.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 824.25
$ws.Range("I94").Value = 692
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 692
$ws.Range("L94").Value = 1750
$ws.Range("M94").Value = -241
$ws.Range("N94").Value = -2652

$ws.Range("H107").Value = 1447.3103
$ws.Range("I107").Value = 1098
$ws.Range("J107").Value = 2223.5557
$ws.Range("K107").Value = 1098
$ws.Range("L107").Value = 2223.5557
$ws.Range("M107").Value = 822
$ws.Range("N107").Value = -6063.5557

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H22").Value = 963.6842
$ws.Range("I22").Value = 314.18182
$ws.Range("J22").Value = 1856.75
$ws.Range("K22").Value = 314.18182
$ws.Range("L22").Value = 1856.75
$ws.Range("M22").Value = 35.81817999999998
$ws.Range("N22").Value = -2556.75

$ws.Range("H31").Value = 1732.77
$ws.Range("I31").Value = 1058.5172
$ws.Range("J31").Value = 2663.8809
$ws.Range("K31").Value = 1058.5172
$ws.Range("L31").Value = 2663.8809
$ws.Range("M31").Value = -763.5172
$ws.Range("N31").Value = -3253.8809

$ws.Range("H34").Value = 1732.77
$ws.Range("I34").Value = 1058.5172
$ws.Range("J34").Value = 2663.8809
$ws.Range("K34").Value = 1058.5172
$ws.Range("L34").Value = 2663.8809
$ws.Range("M34").Value = -856.5172
$ws.Range("N34").Value = -3067.8809

$ws.Range("H99").Value = 2916.3125
$ws.Range("I99").Value = 1624.6
$ws.Range("J99").Value = 3503.4546
$ws.Range("K99").Value = 1624.6
$ws.Range("L99").Value = 3503.4546
$ws.Range("M99").Value = -126.5999999999999
$ws.Range("N99").Value = -6499.4546

$ws.Range("H122").Value = 2318.1562
$ws.Range("I122").Value = 2067.3914
$ws.Range("J122").Value = 2959
$ws.Range("K122").Value = 6202.174199999999
$ws.Range("L122").Value = 8877
$ws.Range("M122").Value = -3752.174199999999
$ws.Range("N122").Value = -13777

$ws.Range("H126").Value = 2916.3125
$ws.Range("I126").Value = 1624.6
$ws.Range("J126").Value = 3503.4546
$ws.Range("K126").Value = 4873.799999999999
$ws.Range("L126").Value = 10510.3638
$ws.Range("M126").Value = -2403.799999999999
$ws.Range("N126").Value = -15450.3638

$ws.Range("H132").Value = 1971.9333
$ws.Range("I132").Value = 1710.2572
$ws.Range("J132").Value = 2887.8
$ws.Range("K132").Value = 5130.7716
$ws.Range("L132").Value = 8663.400000000001
$ws.Range("M132").Value = -2600.7716
$ws.Range("N132").Value = -13723.4

$ws.Range("H134").Value = 12197780
$ws.Range("I134").Value = 16668959
$ws.Range("J134").Value = 3653.9092
$ws.Range("K134").Value = 50006877
$ws.Range("L134").Value = 10961.7276
$ws.Range("M134").Value = -50004342
$ws.Range("N134").Value = -16031.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 739.7857
$ws.Range("J113").Value = 746.85364
$ws.Range("L113").Value = 2240.56092
$ws.Range("N113").Value = -6580.56092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 900.4783
$ws.Range("I22").Value = 429.63635
$ws.Range("J22").Value = 1332.0834
$ws.Range("K22").Value = 429.63635
$ws.Range("L22").Value = 1332.0834
$ws.Range("M22").Value = -134.63635
$ws.Range("N22").Value = -1922.0834

$ws.Range("H27").Value = 900.4783
$ws.Range("I27").Value = 429.63635
$ws.Range("J27").Value = 1332.0834
$ws.Range("K27").Value = 429.63635
$ws.Range("L27").Value = 1332.0834
$ws.Range("M27").Value = -322.63635
$ws.Range("N27").Value = -1546.0834

$ws.Range("H122").Value = 3011.6775
$ws.Range("I122").Value = 2370.9333
$ws.Range("J122").Value = 3612.375
$ws.Range("K122").Value = 7112.7999
$ws.Range("L122").Value = 10837.125
$ws.Range("M122").Value = -4662.7999
$ws.Range("N122").Value = -15737.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 27334.5
$ws.Range("J15").Value = 27334.5
$ws.Range("L15").Value = 27334.5
$ws.Range("N15").Value = -27910.5

$ws.Range("H54").Value = 14012.444
$ws.Range("J54").Value = 14012.444
$ws.Range("L54").Value = 14012.444
$ws.Range("N54").Value = -15052.444

$ws.Range("H81").Value = 2215.6875
$ws.Range("I81").Value = 2000.125
$ws.Range("J81").Value = 2431.25
$ws.Range("K81").Value = 4000.25
$ws.Range("L81").Value = 4862.5
$ws.Range("M81").Value = -2939.25
$ws.Range("N81").Value = -6984.5

$ws.Range("H84").Value = 2215.6875
$ws.Range("I84").Value = 2000.125
$ws.Range("J84").Value = 2431.25
$ws.Range("K84").Value = 20001.25
$ws.Range("L84").Value = 24312.5
$ws.Range("M84").Value = -14697.25
$ws.Range("N84").Value = -34920.5
